$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift data up by one row (drop the old first data row, add new last row data),
# recompute Trad_Prediction (col C) and AI_Prediction (col D) for the new simulation run.

$ws.Cells.Item(2, 1).Value = 45847.60416666666
$ws.Cells.Item(2, 2).Value = 177.8999938964844
$ws.Cells.Item(2, 3).Value = 178.6100006103516
$ws.Cells.Item(2, 4).Value = 179.3526924069895

$ws.Cells.Item(3, 1).Value = 45847.64583333334
$ws.Cells.Item(3, 2).Value = 176.9900054931641
$ws.Cells.Item(3, 3).Value = 177.8999938964844
$ws.Cells.Item(3, 4).Value = 176.8972653981791

$ws.Cells.Item(4, 1).Value = 45847.6875
$ws.Cells.Item(4, 2).Value = 177.2599945068359
$ws.Cells.Item(4, 3).Value = 176.9900054931641
$ws.Cells.Item(4, 4).Value = 178.5697782570394

$ws.Cells.Item(5, 1).Value = 45847.72916666666
$ws.Cells.Item(5, 2).Value = 176.2402038574219
$ws.Cells.Item(5, 3).Value = 177.2599945068359
$ws.Cells.Item(5, 4).Value = 171.2251535048899

$ws.Cells.Item(6, 1).Value = 45847.77083333334
$ws.Cells.Item(6, 2).Value = 176.3600006103516
$ws.Cells.Item(6, 3).Value = 176.2402038574219
$ws.Cells.Item(6, 4).Value = 174.6327813459045

$ws.Cells.Item(7, 1).Value = 45847.8125
$ws.Cells.Item(7, 2).Value = 176.6950073242188
$ws.Cells.Item(7, 3).Value = 176.3600006103516
$ws.Cells.Item(7, 4).Value = 179.5519495055568

$ws.Cells.Item(8, 1).Value = 45848.5625
$ws.Cells.Item(8, 2).Value = 175.1450042724609
$ws.Cells.Item(8, 3).Value = 176.6950073242188
$ws.Cells.Item(8, 4).Value = 174.04584883232

$ws.Cells.Item(9, 1).Value = 45848.60416666666
$ws.Cells.Item(9, 2).Value = 176.3600006103516
$ws.Cells.Item(9, 3).Value = 175.1450042724609
$ws.Cells.Item(9, 4).Value = 176.4885549589004

$ws.Cells.Item(10, 1).Value = 45848.64583333334
$ws.Cells.Item(10, 2).Value = 176.6900024414062
$ws.Cells.Item(10, 3).Value = 176.3600006103516
$ws.Cells.Item(10, 4).Value = 173.6811491612555

$ws.Cells.Item(11, 1).Value = 45848.6875
$ws.Cells.Item(11, 2).Value = 178.1699981689453
$ws.Cells.Item(11, 3).Value = 176.6900024414062
$ws.Cells.Item(11, 4).Value = 181.315162199714

$ws.Cells.Item(12, 1).Value = 45848.72916666666
$ws.Cells.Item(12, 2).Value = 177.6100006103516
$ws.Cells.Item(12, 3).Value = 178.1699981689453
$ws.Cells.Item(12, 4).Value = 178.2389175065274

$ws.Cells.Item(13, 1).Value = 45848.77083333334
$ws.Cells.Item(13, 2).Value = 177.7550048828125
$ws.Cells.Item(13, 3).Value = 177.6100006103516
$ws.Cells.Item(13, 4).Value = 171.1045348304796

$ws.Cells.Item(14, 1).Value = 45848.8125
$ws.Cells.Item(14, 2).Value = 177.6300048828125
$ws.Cells.Item(14, 3).Value = 177.7550048828125
$ws.Cells.Item(14, 4).Value = 178.3230340464171

$ws.Cells.Item(15, 1).Value = 45849.5625
$ws.Cells.Item(15, 2).Value = 177.4450073242188
$ws.Cells.Item(15, 3).Value = 177.6300048828125
$ws.Cells.Item(15, 4).Value = 176.7421679609769

$ws.Cells.Item(16, 1).Value = 45849.60416666666
$ws.Cells.Item(16, 2).Value = 178.9299926757812
$ws.Cells.Item(16, 3).Value = 177.4450073242188
$ws.Cells.Item(16, 4).Value = 179.8917726639792

$ws.Cells.Item(17, 1).Value = 45849.64583333334
$ws.Cells.Item(17, 2).Value = 180.0200042724609
$ws.Cells.Item(17, 3).Value = 178.9299926757812
$ws.Cells.Item(17, 4).Value = 182.9897250093313

$ws.Cells.Item(18, 1).Value = 45849.6875
$ws.Cells.Item(18, 2).Value = 180.5249938964844
$ws.Cells.Item(18, 3).Value = 180.0200042724609
$ws.Cells.Item(18, 4).Value = 180.6756198382009

$ws.Cells.Item(19, 1).Value = 45849.72916666666
$ws.Cells.Item(19, 2).Value = 180.8677978515625
$ws.Cells.Item(19, 3).Value = 180.5249938964844
$ws.Cells.Item(19, 4).Value = 185.8538046365441

$ws.Cells.Item(20, 1).Value = 45849.77083333334
$ws.Cells.Item(20, 2).Value = 180.3863067626953
$ws.Cells.Item(20, 3).Value = 180.8677978515625
$ws.Cells.Item(20, 4).Value = 182.1206646512873

$ws.Cells.Item(21, 1).Value = 45849.8125
$ws.Cells.Item(21, 2).Value = 180.1549987792969
$ws.Cells.Item(21, 3).Value = 180.3863067626953
$ws.Cells.Item(21, 4).Value = 179.3643860467064

$ws.Cells.Item(22, 1).Value = 45852.5625
$ws.Cells.Item(22, 2).Value = 180.1699981689453
$ws.Cells.Item(22, 3).Value = 180.1549987792969
$ws.Cells.Item(22, 4).Value = 179.8528634332398

$ws.Cells.Item(23, 1).Value = 45852.60416666666
$ws.Cells.Item(23, 2).Value = 180.8200073242188
$ws.Cells.Item(23, 3).Value = 180.1699981689453
$ws.Cells.Item(23, 4).Value = 177.3482264237282

$ws.Cells.Item(24, 1).Value = 45852.64583333334
$ws.Cells.Item(24, 2).Value = 180.5850067138672
$ws.Cells.Item(24, 3).Value = 180.8200073242188
$ws.Cells.Item(24, 4).Value = 183.0518243679209

$ws.Cells.Item(25, 1).Value = 45852.6875
$ws.Cells.Item(25, 2).Value = 181.1450042724609
$ws.Cells.Item(25, 3).Value = 180.5850067138672
$ws.Cells.Item(25, 4).Value = 183.8709577186391

$ws.Cells.Item(26, 1).Value = 45852.72916666666
$ws.Cells.Item(26, 2).Value = 181.7299957275391
$ws.Cells.Item(26, 3).Value = 181.1450042724609
$ws.Cells.Item(26, 4).Value = 182.333251961726

$ws.Cells.Item(27, 1).Value = 45852.77083333334
$ws.Cells.Item(27, 2).Value = 181.4799957275391
$ws.Cells.Item(27, 3).Value = 181.7299957275391
$ws.Cells.Item(27, 4).Value = 183.8544686652598

$ws.Cells.Item(28, 1).Value = 45852.8125
$ws.Cells.Item(28, 2).Value = 181.6000061035156
$ws.Cells.Item(28, 3).Value = 181.4799957275391
$ws.Cells.Item(28, 4).Value = 182.1304979635362

$ws.Cells.Item(29, 1).Value = 45853.5625
$ws.Cells.Item(29, 2).Value = 182.5303955078125
$ws.Cells.Item(29, 3).Value = 181.6000061035156
$ws.Cells.Item(29, 4).Value = 181.7619660249487

$ws.Cells.Item(30, 1).Value = 45853.60416666666
$ws.Cells.Item(30, 2).Value = 183.9299926757812
$ws.Cells.Item(30, 3).Value = 182.5303955078125
$ws.Cells.Item(30, 4).Value = 179.6800189838082

$ws.Cells.Item(31, 1).Value = 45853.64583333334
$ws.Cells.Item(31, 2).Value = 183.2350006103516
$ws.Cells.Item(31, 3).Value = 183.9299926757812
$ws.Cells.Item(31, 4).Value = 185.5802657085917

$ws.Cells.Item(32, 1).Value = 45853.6875
$ws.Cells.Item(32, 2).Value = 183.4100036621094
$ws.Cells.Item(32, 3).Value = 183.2350006103516
$ws.Cells.Item(32, 4).Value = 181.2277374065133

$ws.Cells.Item(33, 1).Value = 45853.72916666666
$ws.Cells.Item(33, 2).Value = 183.9149932861328
$ws.Cells.Item(33, 3).Value = 183.4100036621094
$ws.Cells.Item(33, 4).Value = 179.3471344066051

$ws.Cells.Item(34, 1).Value = 45853.77083333334
$ws.Cells.Item(34, 2).Value = 183.2100067138672
$ws.Cells.Item(34, 3).Value = 183.9149932861328
$ws.Cells.Item(34, 4).Value = 188.8395410456848

$ws.Cells.Item(35, 1).Value = 45853.8125
$ws.Cells.Item(35, 2).Value = 181.9600067138672
$ws.Cells.Item(35, 3).Value = 183.2100067138672
$ws.Cells.Item(35, 4).Value = 180.9906202541765

# Remove the now-obsolete last row (old row 36), shrinking the sheet from A1:D36 to A1:D35
$ws.Rows.Item(36).Delete()

Write-Host "Applied simulation data shift + tradicional agent recompute"